$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 29.20950566666667
$ws.Range("H2").Value = 87.628517
$ws.Range("I2").Value = 0.01829497698069002
$ws.Range("J2").Value = 0.01840828041918582
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1744923333333333
$ws.Range("N2").Value = 0.523477
$ws.Range("O2").Value = 0.07822917822503123
$ws.Range("P2").Value = 0.08239975633156223
$ws.Range("Q2").Value = 5.096834799289889
$ws.Range("R2").Value = 45.871513193609
$ws.Range("S2").Value = 0.001431201014845243
$ws.Range("T2").Value = 0.00151683782102398

$ws.Range("G3").Value = 29.20950566666667
$ws.Range("H3").Value = 87.628517
$ws.Range("I3").Value = 0.01829497698069002
$ws.Range("J3").Value = 0.01840828041918582
$ws.Range("O3").Value = 0.5425629886502931
$ws.Range("P3").Value = 0.5714882742434749
$ws.Range("Q3").Value = 35.34939244030944
$ws.Range("R3").Value = 318.144531962785
$ws.Range("S3").Value = 0.009926177387931492
$ws.Range("T3").Value = 0.01052011640855046

$ws.Range("G4").Value = 29.20950566666667
$ws.Range("H4").Value = 87.628517
$ws.Range("I4").Value = 0.01829497698069002
$ws.Range("J4").Value = 0.01840828041918582
$ws.Range("M4").Value = 0.1427166666666667
$ws.Range("N4").Value = 0.42815
$ws.Range("O4").Value = 0.06398337015197826
$ws.Range("P4").Value = 0.06739447133944447
$ws.Range("Q4").Value = 4.168683283727779
$ws.Range("R4").Value = 37.51814955355
$ws.Range("S4").Value = 0.001170574284077411
$ws.Range("T4").Value = 0.001240616327119276

$ws.Range("G5").Value = 29.20950566666667
$ws.Range("H5").Value = 87.628517
$ws.Range("I5").Value = 0.01829497698069002
$ws.Range("J5").Value = 0.01840828041918582
$ws.Range("M5").Value = 0.3386875
$ws.Range("N5").Value = 0.6773750000000001
$ws.Range("O5").Value = 0.15184188493529
$ws.Range("P5").Value = 0.1066246175956001
$ws.Range("Q5").Value = 9.892894450479169
$ws.Range("R5").Value = 59.35736670287501
$ws.Range("S5").Value = 0.002777943789595713
$ws.Range("T5").Value = 0.001962775860288262

$ws.Range("G6").Value = 29.20950566666667
$ws.Range("H6").Value = 87.628517
$ws.Range("I6").Value = 0.01829497698069002
$ws.Range("J6").Value = 0.01840828041918582
$ws.Range("M6").Value = 0.3644293333333333
$ws.Range("N6").Value = 1.093288
$ws.Range("O6").Value = 0.1633825780374074
$ws.Range("P6").Value = 0.1720928804899184
$ws.Range("Q6").Value = 10.64480067709956
$ws.Range("R6").Value = 95.80320609389601
$ws.Range("S6").Value = 0.00298908050424016
$ws.Range("T6").Value = 0.00316793400220385

$ws.Range("I7").Value = 0.913374480506715
$ws.Range("J7").Value = 0.9190311407684336
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.1744923333333333
$ws.Range("N7").Value = 0.523477
$ws.Range("O7").Value = 0.07822917822503123
$ws.Range("P7").Value = 0.08239975633156223
$ws.Range("Q7").Value = 254.4588518446099
$ws.Range("R7").Value = 2290.129666601489
$ws.Range("S7").Value = 0.07145253502175512
$ws.Range("T7").Value = 0.07572794206043659

$ws.Range("I8").Value = 0.913374480506715
$ws.Range("J8").Value = 0.9190311407684336
$ws.Range("O8").Value = 0.5425629886502931
$ws.Range("P8").Value = 0.5714882742434749
$ws.Range("S8").Value = 0.4955631879006321
$ws.Range("T8").Value = 0.5252155206137642

$ws.Range("I9").Value = 0.913374480506715
$ws.Range("J9").Value = 0.9190311407684336
$ws.Range("M9").Value = 0.1427166666666667
$ws.Range("N9").Value = 0.42815
$ws.Range("O9").Value = 0.06398337015197826
$ws.Range("P9").Value = 0.06739447133944447
$ws.Range("Q9").Value = 208.1210013377278
$ws.Range("R9").Value = 1873.08901203955
$ws.Range("S9").Value = 0.058440777473632
$ws.Range("T9").Value = 0.06193761787657515

$ws.Range("I10").Value = 0.913374480506715
$ws.Range("J10").Value = 0.9190311407684336
$ws.Range("M10").Value = 0.3386875
$ws.Range("N10").Value = 0.6773750000000001
$ws.Range("O10").Value = 0.15184188493529
$ws.Range("P10").Value = 0.1066246175956001
$ws.Range("Q10").Value = 493.9015413329791
$ws.Range("R10").Value = 2963.409247997875
$ws.Range("S10").Value = 0.1386885027719309
$ws.Range("T10").Value = 0.09799134394288241

$ws.Range("I11").Value = 0.913374480506715
$ws.Range("J11").Value = 0.9190311407684336
$ws.Range("M11").Value = 0.3644293333333333
$ws.Range("N11").Value = 1.093288
$ws.Range("O11").Value = 0.1633825780374074
$ws.Range("P11").Value = 0.1720928804899184
$ws.Range("Q11").Value = 531.4403674191796
$ws.Range("R11").Value = 4782.963306772616
$ws.Range("S11").Value = 0.1492294773387648
$ws.Range("T11").Value = 0.1581587162747754

$ws.Range("G12").Value = 57.98602933333333
$ws.Range("H12").Value = 173.958088
$ws.Range("I12").Value = 0.03631876156896331
$ws.Range("J12").Value = 0.03654368891224535
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.1744923333333333
$ws.Range("N12").Value = 0.523477
$ws.Range("O12").Value = 0.07822917822503123
$ws.Range("P12").Value = 0.08239975633156223
$ws.Range("Q12").Value = 10.11811755910844
$ws.Range("R12").Value = 91.06305803197598
$ws.Range("S12").Value = 0.002841186871690846
$ws.Range("T12").Value = 0.003011191061825429

$ws.Range("G13").Value = 57.98602933333333
$ws.Range("H13").Value = 173.958088
$ws.Range("I13").Value = 0.03631876156896331
$ws.Range("J13").Value = 0.03654368891224535
$ws.Range("O13").Value = 0.5425629886502931
$ws.Range("P13").Value = 0.5714882742434749
$ws.Range("Q13").Value = 70.17478934258222
$ws.Range("R13").Value = 631.57310408324
$ws.Range("S13").Value = 0.01970521582093414
$ws.Range("T13").Value = 0.0208842897109495

$ws.Range("G14").Value = 57.98602933333333
$ws.Range("H14").Value = 173.958088
$ws.Range("I14").Value = 0.03631876156896331
$ws.Range("J14").Value = 0.03654368891224535
$ws.Range("M14").Value = 0.1427166666666667
$ws.Range("N14").Value = 0.42815
$ws.Range("O14").Value = 0.06398337015197826
$ws.Range("P14").Value = 0.06739447133944447
$ws.Range("Q14").Value = 8.27557281968889
$ws.Range("R14").Value = 74.48015537719999
$ws.Range("S14").Value = 0.002323796764928422
$ws.Range("T14").Value = 0.002462842595033893

$ws.Range("G15").Value = 57.98602933333333
$ws.Range("H15").Value = 173.958088
$ws.Range("I15").Value = 0.03631876156896331
$ws.Range("J15").Value = 0.03654368891224535
$ws.Range("M15").Value = 0.3386875
$ws.Range("N15").Value = 0.6773750000000001
$ws.Range("O15").Value = 0.15184188493529
$ws.Range("P15").Value = 0.1066246175956001
$ws.Range("Q15").Value = 19.63914330983333
$ws.Range("R15").Value = 117.834859859
$ws.Range("S15").Value = 0.005514709215146759
$ws.Range("T15").Value = 0.003896456855800733

$ws.Range("G16").Value = 57.98602933333333
$ws.Range("H16").Value = 173.958088
$ws.Range("I16").Value = 0.03631876156896331
$ws.Range("J16").Value = 0.03654368891224535
$ws.Range("M16").Value = 0.3644293333333333
$ws.Range("N16").Value = 1.093288
$ws.Range("O16").Value = 0.1633825780374074
$ws.Range("P16").Value = 0.1720928804899184
$ws.Range("Q16").Value = 21.13181001259377
$ws.Range("R16").Value = 190.186290113344
$ws.Range("S16").Value = 0.005933852896263143
$ws.Range("T16").Value = 0.006288908688635794

$ws.Range("G17").Value = 29.481085
$ws.Range("H17").Value = 58.96217
$ws.Range("I17").Value = 0.01846507700595112
$ws.Range("J17").Value = 0.01238628926567028
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.1744923333333333
$ws.Range("N17").Value = 0.523477
$ws.Range("O17").Value = 0.07822917822503123
$ws.Range("P17").Value = 0.08239975633156223
$ws.Range("Q17").Value = 5.144223310848333
$ws.Range("R17").Value = 30.86533986509
$ws.Range("S17").Value = 0.001444507800037477
$ws.Range("T17").Value = 0.001020627217343476

$ws.Range("G18").Value = 29.481085
$ws.Range("H18").Value = 58.96217
$ws.Range("I18").Value = 0.01846507700595112
$ws.Range("J18").Value = 0.01238628926567028
$ws.Range("O18").Value = 0.5425629886502931
$ws.Range("P18").Value = 0.5714882742434749
$ws.Range("Q18").Value = 35.67805820214166
$ws.Range("R18").Value = 214.06834921285
$ws.Range("S18").Value = 0.01001846736600665
$ws.Range("T18").Value = 0.007078619076718386

$ws.Range("G19").Value = 29.481085
$ws.Range("H19").Value = 58.96217
$ws.Range("I19").Value = 0.01846507700595112
$ws.Range("J19").Value = 0.01238628926567028
$ws.Range("M19").Value = 0.1427166666666667
$ws.Range("N19").Value = 0.42815
$ws.Range("O19").Value = 0.06398337015197826
$ws.Range("P19").Value = 0.06739447133944447
$ws.Range("Q19").Value = 4.207442180916667
$ws.Range("R19").Value = 25.2446530855
$ws.Range("S19").Value = 0.001181457856956553
$ws.Range("T19").Value = 0.0008347674169172842

$ws.Range("G20").Value = 29.481085
$ws.Range("H20").Value = 58.96217
$ws.Range("I20").Value = 0.01846507700595112
$ws.Range("J20").Value = 0.01238628926567028
$ws.Range("M20").Value = 0.3386875
$ws.Range("N20").Value = 0.6773750000000001
$ws.Range("O20").Value = 0.15184188493529
$ws.Range("P20").Value = 0.1066246175956001
$ws.Range("Q20").Value = 9.984874975937501
$ws.Range("R20").Value = 39.93949990375
$ws.Range("S20").Value = 0.0028037720980589
$ws.Range("T20").Value = 0.00132068335638058

$ws.Range("G21").Value = 29.481085
$ws.Range("H21").Value = 58.96217
$ws.Range("I21").Value = 0.01846507700595112
$ws.Range("J21").Value = 0.01238628926567028
$ws.Range("M21").Value = 0.3644293333333333
$ws.Range("N21").Value = 1.093288
$ws.Range("O21").Value = 0.1633825780374074
$ws.Range("P21").Value = 0.1720928804899184
$ws.Range("Q21").Value = 10.74377215249333
$ws.Range("R21").Value = 64.46263291496
$ws.Range("S21").Value = 0.003016871884891547
$ws.Range("T21").Value = 0.002131592198310554

$ws.Range("G22").Value = 21.628479
$ws.Range("H22").Value = 64.885437
$ws.Range("I22").Value = 0.01354670393768061
$ws.Range("J22").Value = 0.01363060063446486
$ws.Range("K22").Value = 2
$ws.Range("L22").Value = 0.6666666666666666
$ws.Range("M22").Value = 0.1744923333333333
$ws.Range("N22").Value = 0.523477
$ws.Range("O22").Value = 0.07822917822503123
$ws.Range("P22").Value = 0.08239975633156223
$ws.Range("Q22").Value = 3.774003767161
$ws.Range("R22").Value = 33.966033904449
$ws.Range("S22").Value = 0.001059747516702549
$ws.Range("T22").Value = 0.001123158170932742

$ws.Range("G23").Value = 21.628479
$ws.Range("H23").Value = 64.885437
$ws.Range("I23").Value = 0.01354670393768061
$ws.Range("J23").Value = 0.01363060063446486
$ws.Range("O23").Value = 0.5425629886502931
$ws.Range("P23").Value = 0.5714882742434749
$ws.Range("Q23").Value = 26.174821333265
$ws.Range("R23").Value = 235.573391999385
$ws.Range("S23").Value = 0.007349940174788685
$ws.Range("T23").Value = 0.007789728433492338

$ws.Range("G24").Value = 21.628479
$ws.Range("H24").Value = 64.885437
$ws.Range("I24").Value = 0.01354670393768061
$ws.Range("J24").Value = 0.01363060063446486
$ws.Range("M24").Value = 0.1427166666666667
$ws.Range("N24").Value = 0.42815
$ws.Range("O24").Value = 0.06398337015197826
$ws.Range("P24").Value = 0.06739447133944447
$ws.Range("Q24").Value = 3.08674442795
$ws.Range("R24").Value = 27.78069985155
$ws.Range("S24").Value = 0.0008667637723838799
$ws.Range("T24").Value = 0.0009186271237988558

$ws.Range("G25").Value = 21.628479
$ws.Range("H25").Value = 64.885437
$ws.Range("I25").Value = 0.01354670393768061
$ws.Range("J25").Value = 0.01363060063446486
$ws.Range("M25").Value = 0.3386875
$ws.Range("N25").Value = 0.6773750000000001
$ws.Range("O25").Value = 0.15184188493529
$ws.Range("P25").Value = 0.1066246175956001
$ws.Range("Q25").Value = 7.325295481312501
$ws.Range("R25").Value = 43.951772887875
$ws.Range("S25").Value = 0.002056957060557739
$ws.Range("T25").Value = 0.001453357580248161

$ws.Range("G26").Value = 21.628479
$ws.Range("H26").Value = 64.885437
$ws.Range("I26").Value = 0.01354670393768061
$ws.Range("J26").Value = 0.01363060063446486
$ws.Range("M26").Value = 0.3644293333333333
$ws.Range("N26").Value = 1.093288
$ws.Range("O26").Value = 0.1633825780374074
$ws.Range("P26").Value = 0.1720928804899184
$ws.Range("Q26").Value = 7.882052182983999
$ws.Range("R26").Value = 70.93846964685599
$ws.Range("S26").Value = 0.002213295413247757
$ws.Range("T26").Value = 0.002345729325992767
